# Atualiza dados BIBI: remove vendas atipicas de 2025-06-11 e corrige
# os valores de estoque_atualizado para os registros remanescentes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove as linhas referentes as vendas de 2025-06-11 (id_venda 357349, 357392, 357402)
$ws.Rows("2:4").Delete() | Out-Null

# Apos a remocao, as linhas restantes foram deslocadas para cima (linhas 2-10).
# Corrige a coluna G (estoque_atualizado) para os pedidos impactados.
$ws.Range("G3").Value = -293     # id_venda 358540
$ws.Range("G6").Value = -439     # id_venda 362396
$ws.Range("G7").Value = -293     # id_venda 362404
$ws.Range("G8").Value = -439     # id_venda 365782
$ws.Range("G9").Value = -439     # id_venda 366707
